# Update the build timestamp embedded in the workbook text from
# "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")

$a2 = $aboutSheet.Range("A2").Value2
$aboutSheet.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $aboutSheet.Range("A6").Value2
$aboutSheet.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 7; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # Column S
    $val = $cell.Value2
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
